$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Update filmgold row (row 2) values: k and d recalculated
$ws.Range("B2").Value = 140.263
$ws.Range("C2").Value = 0.0145

# Insert 5 new rows after row 2 (before the old pyrex row), shifting
# pyrex/gold/quartz/seed rows down to make room for new materials
$ws.Range("A3:D7").Insert()

# Fill in the new material rows (BK7, Borofloat, CaF2, ZnSe, ZnS).
# New shared-string entries are created in this order: BK7, Borofloat,
# CaF2, ZnSe, ZnS - so set A4 before A3 to reproduce that ordering.
$ws.Range("A4").Value = "BK7"
$ws.Range("B4").Value = 1.114
$ws.Range("C4").Value = 0.000000516676
$ws.Range("D4").Value = 2510

$ws.Range("A3").Value = "Borofloat"
$ws.Range("B3").Value = 1.2
$ws.Range("C3").Value = 0.000000657
$ws.Range("D3").Value = 2200

$ws.Range("A5").Value = "CaF2"
$ws.Range("B5").Value = 9.71
$ws.Range("C5").Value = 0.0000035755
$ws.Range("D5").Value = 3180

$ws.Range("A6").Value = "ZnSe"
$ws.Range("B6").Value = 18
$ws.Range("C6").Value = 0.000010075
$ws.Range("D6").Value = 5270

$ws.Range("A7").Value = "ZnS"
$ws.Range("B7").Value = 27.2
$ws.Range("C7").Value = 0.0000129133
$ws.Range("D7").Value = 4090

# Update selection to match the recorded cursor position after edits
$ws.Activate()
$ws.Range("D18").Select()
